$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part A: split the "(5 points)" validation-user-story paragraph so the
# hidden "_GoBack" bookmark moves from before the "(10 points): As a user..."
# paragraph to sit inside "provide" (between "provid" and "e invalid input.").
# ---------------------------------------------------------------------------

# The _GoBack bookmark currently sits just before the "(10 points): As a
# user, I want to be able to search..." paragraph. Remove it from there.
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# Locate the point right after "provid" (i.e. before "e invalid input.").
$find1 = $d.Content
[void]$find1.Find.Execute("As a developer, I want to run validation on any user input, ensuring that a user is re-prompted when they provid")
$splitPoint = $find1.End

# Re-adding the bookmark here (a collapsed, zero-length range) both places it
# at the correct spot and forces the run to split cleanly into two runs.
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange)

# ---------------------------------------------------------------------------
# Part B: highlight the "(10 points): ... highlighted in yellow." paragraph
# in yellow (matching the other user-story paragraphs), while leaving the
# single trailing space before the paragraph mark un-highlighted.
# ---------------------------------------------------------------------------

# Highlight the bold "(10 points):" run.
$f1 = $d.Content.Find
$f1.Text = "(10 points):"
$f1.Replacement.Text = "(10 points):"
$f1.Replacement.Highlight = $true
[void]$f1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 1)

# Highlight the remaining sentence (leading space included, trailing space
# excluded so it stays un-highlighted in its own run).
$f2 = $d.Content.Find
$f2.Text = " As a developer, I will send a copy of these user stories to the instructors at the end of each workday, with user stories fully implemented highlighted in green and partially implemented highlighted in yellow."
$f2.Replacement.Text = " As a developer, I will send a copy of these user stories to the instructors at the end of each workday, with user stories fully implemented highlighted in green and partially implemented highlighted in yellow."
$f2.Replacement.Highlight = $true
[void]$f2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 1)
